# Actualización desde MV -datos-
# Appends 30 new daily rows (05-08-2021 .. 03-09-2021) to the bottom of the
# "Diaria" sheet, continuing the existing Serie/UF/Pesos table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$startRow = 191
$dates = @(
    "05-08-2021","06-08-2021","07-08-2021","08-08-2021","09-08-2021",
    "10-08-2021","11-08-2021","12-08-2021","13-08-2021","14-08-2021",
    "15-08-2021","16-08-2021","17-08-2021","18-08-2021","19-08-2021",
    "20-08-2021","21-08-2021","22-08-2021","23-08-2021","24-08-2021",
    "25-08-2021","26-08-2021","27-08-2021","28-08-2021","29-08-2021",
    "30-08-2021","31-08-2021","01-09-2021","02-09-2021","03-09-2021"
)

for ($i = 0; $i -lt $dates.Length; $i++) {
    $r = $startRow + $i

    # Write the date as a formula returning a text string, then flatten it
    # back to a plain value via copy / paste-special-values. This keeps the
    # cell a genuine shared-string text cell (matching the other date rows)
    # instead of letting Excel's smart-entry re-interpret "dd-mm-yyyy" text
    # as a date serial number (which would also mint a brand-new style).
    $cellA = $ws.Range("A$r")
    $cellA.Formula = '="' + $dates[$i] + '"'
    $cellA.Copy()
    $cellA.PasteSpecial(-4163)

    $ws.Range("B$r").Value = 449
    $ws.Range("C$r").Value = 0
}
